$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-7): dates become text, paths Shapes_N -> Group_N ---

# Row 2: 20210413, Shapes_1 -> Group_1, date 44299 -> "13/04/2021" (text)
$ws.Range("A2").Value = "13/04/2021"
$ws.Range("B2").Value = "F:\PhD, PMMH, ESPCI\Processing\20210413-Actin\results\Group_1"

# Row 3: 20210430, Shapes_1 -> Group_1, date 44316 -> "30/04/2021" (text)
$ws.Range("A3").Value = "30/04/2021"
$ws.Range("B3").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Group_1"

# Row 4: 20210430, Shapes_2 -> Group_2, date 44316 -> "30/04/2021" (text)
$ws.Range("A4").Value = "30/04/2021"
$ws.Range("B4").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Group_2"

# Row 5: 20200925, Shapes_1 -> Group_1, date 44099 -> "25/09/2020" (text)
$ws.Range("A5").Value = "25/09/2020"
$ws.Range("B5").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Group_1"
$ws.Range("C5").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\circlesforPAs1.mat"
$ws.Range("D5").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Figures"

# Row 6: 20200925, Shapes_2 -> Group_2, date 44099 -> "25/09/2020" (text)
$ws.Range("A6").Value = "25/09/2020"
$ws.Range("B6").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Group_2"
$ws.Range("C6").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\circlesforPAs2.mat"
$ws.Range("D6").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Figures"

# Row 7: 20201001, Shapes_1 -> Group_1, date 44105 -> 43840 (stays numeric)
$ws.Range("A7").Value = 43840
$ws.Range("B7").Value = "F:\PhD, PMMH, ESPCI\Processing\20201001-Actin\results\Group_1"
$ws.Range("C7").Value = "F:\PhD, PMMH, ESPCI\Processing\20201001-Actin\results\circlesforPAs1.mat"
$ws.Range("D7").Value = "F:\PhD, PMMH, ESPCI\Processing\20201001-Actin\results\Figures"

# --- Add new rows (8-10): 20211029-Actin experiment, 3 groups ---

$ws.Range("A8").Value = "29/10/2021"
$ws.Range("B8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_1"
$ws.Range("C8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D8").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 400
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Formula = "=G8/F8/E8"
$ws.Range("I8").Value = 0.1
$ws.Range("J8").Value = 20

$ws.Range("A9").Value = "29/10/2021"
$ws.Range("B9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_2"
$ws.Range("C9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D9").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E9").Value = 55
$ws.Range("F9").Value = 400
$ws.Range("G9").Value = 1
$ws.Range("H9").Formula = "=G9/F9/E9"
$ws.Range("I9").Value = 0.1
$ws.Range("J9").Value = 20

$ws.Range("A10").Value = "29/10/2021"
$ws.Range("B10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Group_3"
$ws.Range("C10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\circlesforPAs20211029_S15.mat"
$ws.Range("D10").Value = "G:\PhD, PMMH, ESPCI\Processing\20211029-Actin\results\Figures"
$ws.Range("E10").Value = 55
$ws.Range("F10").Value = 400
$ws.Range("G10").Value = 1.5
$ws.Range("H10").Formula = "=G10/F10/E10"
$ws.Range("I10").Value = 0.1
$ws.Range("J10").Value = 20

# --- Column width / selection cosmetic updates ---
# COM ColumnWidth vs. stored XML width differ by 5/6 of a character unit on this
# engine; 81+1/6 round-trips to an XML width of exactly 82 (the target value).
$ws.Columns.Item(3).ColumnWidth = 81.16666666666667
$ws.Range("E13").Select()
